$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "first_name"/"last_name" columns and repurpose the space: D
# becomes a new "phone" column, E becomes "is_admin" (moved in from the
# old, now-removed, G column).
$ws.Range("D1").Value = "phone"
$ws.Range("E1").Value = "is_admin"

$ws.Range("D2").Value = 12345678
$ws.Range("E2").Value = 1

$ws.Range("D3").Value = 77776666
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 44447777
$ws.Range("E4").Value = 0

# The old "dob" (F) and "is_admin" (G) columns are no longer part of the
# table at all - wipe them out completely (values + formatting) rather
# than leaving blanks behind.
$ws.Range("F1:G4").Clear() | Out-Null

# A couple of the seeded passwords changed.
$ws.Range("B3").Value = 19841984
$ws.Range("B4").Value = 19911991

# The leading "username" column now gets an explicit width, and the
# worksheet's remembered selection moved too.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("C15").Select() | Out-Null
